$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix font family for the "LoadingUI" rich-text cells (H2:H5 all share the
# "UI/ChronoBlade_forest_wallpaper" string) - set the font Family to 3
# (matches the "宋体" / Chinese font family classification used elsewhere).
$ws.Range("H2:H5").Font.Family = 3

# Correct the previously-wrong CanClone flag for row 6 (DemoWinter scene)
$ws.Range("M6").Value = 1

# Move the active selection as last left by the editor
$ws.Range("N12").Select()
